# Arreglos en el informe de recategorizacion
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 2: "Montos para la nueva categoria"
$ws.Range("A2").Value2 = "Montos para la nueva categoria"

# Update existing value: Ingresos Brutos Devengados (B5) from 723270.02 to 725070.02
$ws.Range("B5").Value2 = 725070.02

# New rows 8-10: facturacion remaining detail
$ws.Range("A8").Value2 = "Hasta el 30 jun 2020 todavia podes facturar"
$ws.Range("B8").Value2 = 109886.98

$ws.Range("A9").Value2 = "Faltan 3 meses para la proxima recategorizacion"
$ws.Range("B9").Value2 = 3

$ws.Range("A10").Value2 = "Por mes podes facturar"
$ws.Range("B10").Formula = "=+B8/B9"

# New rows 13-15: cuota / aportes / total
$ws.Range("A13").Value2 = "Cuota Mensual"
$ws.Range("B13").Value2 = 1739.48

$ws.Range("A14").Value2 = "Aporte autonomo / Obra social"
$ws.Range("B14").Value2 = 5256.36

$ws.Range("A15").Value2 = "Total"
$ws.Range("B15").Value2 = 6995.84

# Move the active selection to A12, as in the edited workbook
$ws.Range("A12").Select()
